$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 600
$ws.Range("I28").Value = 436.36365
$ws.Range("J28").Value = 1200
$ws.Range("K28").Value = 436.36365
$ws.Range("L28").Value = 1200
$ws.Range("M28").Value = 48.63634999999999
$ws.Range("N28").Value = -2170

# Row 113
$ws.Range("H113").Value = 3716.6667
$ws.Range("I113").Value = 2799.4443
$ws.Range("J113").Value = 5092.5
$ws.Range("K113").Value = 2799.4443
$ws.Range("L113").Value = 5092.5
$ws.Range("M113").Value = 454.5556999999999
$ws.Range("N113").Value = -11600.5

# Row 129
$ws.Range("H129").Value = 1175.2533
$ws.Range("I129").Value = 791.55554
$ws.Range("J129").Value = 1227.5758
$ws.Range("K129").Value = 2374.66662
$ws.Range("L129").Value = 3682.7274
$ws.Range("M129").Value = 2625.33338
$ws.Range("N129").Value = -13682.7274

# Row 132
$ws.Range("H132").Value = 7187.825
$ws.Range("I132").Value = 6133.6294
$ws.Range("J132").Value = 9377.308000000001
$ws.Range("K132").Value = 18400.8882
$ws.Range("L132").Value = 28131.924
$ws.Range("M132").Value = -15870.8882
$ws.Range("N132").Value = -33191.924

$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 16
$ws.Range("H16").Value = 700
$ws.Range("I16").Value = 700
$ws.Range("K16").Value = 700
$ws.Range("M16").Value = -413

# Row 61
$ws.Range("H61").Value = 4915.0835
$ws.Range("I61").Value = 7311.385
$ws.Range("J61").Value = 2083.0908
$ws.Range("K61").Value = 7311.385
$ws.Range("L61").Value = 2083.0908
$ws.Range("M61").Value = -7099.385
$ws.Range("N61").Value = -2507.0908

# Row 74
$ws.Range("H74").Value = 1212.9512
$ws.Range("I74").Value = 1008.9643
$ws.Range("J74").Value = 1652.3077
$ws.Range("K74").Value = 1008.9643
$ws.Range("L74").Value = 1652.3077
$ws.Range("M74").Value = -134.9643
$ws.Range("N74").Value = -3400.3077

# Row 77
$ws.Range("H77").Value = 1212.9512
$ws.Range("I77").Value = 1008.9643
$ws.Range("J77").Value = 1652.3077
$ws.Range("K77").Value = 5044.8215
$ws.Range("L77").Value = 8261.538500000001
$ws.Range("M77").Value = -676.8215
$ws.Range("N77").Value = -16997.5385

# Row 136
$ws.Range("H136").Value = 4915.0835
$ws.Range("I136").Value = 7311.385
$ws.Range("J136").Value = 2083.0908
$ws.Range("K136").Value = 21934.155
$ws.Range("L136").Value = 6249.2724
$ws.Range("M136").Value = -19384.155
$ws.Range("N136").Value = -11349.2724

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 500
$ws.Range("K8").Value = 500
$ws.Range("M8").Value = -360

# Row 10
$ws.Range("H10").Value = 28001.666
$ws.Range("I10").Value = 1999.5
$ws.Range("K10").Value = 1999.5
$ws.Range("M10").Value = -1859.5

# Row 22
$ws.Range("H22").Value = 1557.4615
$ws.Range("I22").Value = 1877.4445
$ws.Range("K22").Value = 1877.4445
$ws.Range("M22").Value = -1704.4445

# Row 134
$ws.Range("H134").Value = 5142.0444
$ws.Range("I134").Value = 1852.6923
$ws.Range("K134").Value = 5558.0769
$ws.Range("M134").Value = -3023.0769

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 7336.7646
$ws.Range("I58").Value = 3079
$ws.Range("J58").Value = 12126.75
$ws.Range("K58").Value = 3079
$ws.Range("L58").Value = 12126.75
$ws.Range("M58").Value = -2876
$ws.Range("N58").Value = -12532.75

# Row 122
$ws.Range("H122").Value = 333336930
$ws.Range("J122").Value = 10000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -34900

# Row 136
$ws.Range("H136").Value = 7336.7646
$ws.Range("I136").Value = 3079
$ws.Range("J136").Value = 12126.75
$ws.Range("K136").Value = 9237
$ws.Range("L136").Value = 36380.25
$ws.Range("M136").Value = -6687
$ws.Range("N136").Value = -41480.25

$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 816.8095
$ws.Range("I92").Value = 513.7273
$ws.Range("J92").Value = 1150.2
$ws.Range("K92").Value = 1541.1819
$ws.Range("L92").Value = 3450.6
$ws.Range("M92").Value = -293.1819
$ws.Range("N92").Value = -5946.6

# Row 107
$ws.Range("H107").Value = 38462010
$ws.Range("I107").Value = 125000410
$ws.Range("J107").Value = 500.33334
$ws.Range("K107").Value = 375001230
$ws.Range("L107").Value = 1501.00002
$ws.Range("M107").Value = -374999310
$ws.Range("N107").Value = -5341.000019999999

# Row 113
$ws.Range("H113").Value = 481.0606
$ws.Range("I113").Value = 412.42856
$ws.Range("J113").Value = 601.1667
$ws.Range("K113").Value = 1237.28568
$ws.Range("L113").Value = 1803.5001
$ws.Range("M113").Value = 932.71432
$ws.Range("N113").Value = -6143.5001

# Row 131
$ws.Range("H131").Value = 1031.5588
$ws.Range("I131").Value = 554.44446
$ws.Range("J131").Value = 1203.32
$ws.Range("K131").Value = 1663.33338
$ws.Range("L131").Value = 3609.96
$ws.Range("M131").Value = 3376.66662
$ws.Range("N131").Value = -13689.96

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 6018.875
$ws.Range("I122").Value = 5525.1665
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 16575.4995
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -14125.4995
$ws.Range("N122").Value = -27400

# Row 126
$ws.Range("H126").Value = 10419126
$ws.Range("I126").Value = 19232376
$ws.Range("J126").Value = 3466.7273
$ws.Range("K126").Value = 57697128
$ws.Range("L126").Value = 10400.1819
$ws.Range("M126").Value = -57694658
$ws.Range("N126").Value = -15340.1819

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 32559.656
$ws.Range("I132").Value = 40651.52
$ws.Range("J132").Value = 5249.625
$ws.Range("K132").Value = 121954.56
$ws.Range("L132").Value = 15748.875
$ws.Range("M132").Value = -119424.56
$ws.Range("N132").Value = -20808.875

# Row 136
$ws.Range("H136").Value = 2676.4707
$ws.Range("I136").Value = 2433.3333
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 7299.999899999999
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = -4749.999899999999
$ws.Range("N136").Value = -13950

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 962.8
$ws.Range("I113").Value = 407.7143
$ws.Range("J113").Value = 1448.5
$ws.Range("K113").Value = 1223.1429
$ws.Range("L113").Value = 4345.5
$ws.Range("M113").Value = 946.8571000000002
$ws.Range("N113").Value = -8685.5

# Row 122
$ws.Range("H122").Value = 2641
$ws.Range("I122").Value = 1811.7142
$ws.Range("J122").Value = 3470.2856
$ws.Range("K122").Value = 5435.142599999999
$ws.Range("L122").Value = 10410.8568
$ws.Range("M122").Value = -2985.142599999999
$ws.Range("N122").Value = -15310.8568

Write-Output "Applied Mandragora_Profits updates"